# Repurpose geothermal as pumped hydro
#
# 1. On the PTCF sheet, the "geothermal" row's Summer peak time capacity
#    factor (B10) used to mirror the lignite value (=B13, 0.8). Since the
#    India EPS repurposes the geothermal plant type to represent pumped
#    hydro capacity, set it to a literal capacity factor of 1.
# 2. Add an explanatory note about this repurposing to the About sheet.

$wb = $excel.ActiveWorkbook

$ptcf = $wb.Worksheets.Item("PTCF")
$ptcf.Range("B10").Value = 1

$about = $wb.Worksheets.Item("About")
$note = $about.Range("A29")
$note.Value = "In the India EPS, the geothermal plant type is repurposed as pumped hydro capacity."
$note.Font.Color = 0
$note.VerticalAlignment = -4108
